$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

# Row 47 - Spillkråka (Dryocopus martius)
$ws.Range("A47").Value = 131082539
$ws.Range("B47").Value = 57881
$ws.Range("D47").Value = "NT"
$ws.Range("E47").Value = 100049
$ws.Range("F47").Value = "Spillkråka"
$ws.Range("G47").Value = "Dryocopus martius"
$ws.Range("H47").Value = "(Linnaeus, 1758)"
# Leading apostrophe forces text storage for values that look numeric/date-like
$ws.Range("I47").Value = "'1"
$ws.Range("K47").Value = "adult"
$ws.Range("L47").Value = "hane"
$ws.Range("M47").Value = "permanent revir"
# A lone apostrophe round-trips as an explicit empty text cell (unlike
# assigning "" directly, which leaves the cell truly blank/absent).
$ws.Range("N47").Value = "'"
$ws.Range("P47").Value = "Glupen- nedklassad naturvårdsskog, Ög"
$ws.Range("Q47").Value = 571452
$ws.Range("R47").Value = 6465778
$ws.Range("S47").Value = 10
$ws.Range("T47").Value = "Östergötland"
$ws.Range("U47").Value = "Söderköping"
$ws.Range("V47").Value = "Östergötland"
$ws.Range("W47").Value = "Östra Ryd"
$ws.Range("Y47").Value = "'2026-02-07"
$ws.Range("AA47").Value = "'2026-02-07"
$ws.Range("AD47").Value = $false
$ws.Range("AE47").Value = $false
$ws.Range("AG47").Value = $false
$ws.Range("AT47").Value = "'"
$ws.Range("AW47").Value = "Steve Daurer"
$ws.Range("AX47").Value = "Steve Daurer"
$ws.Range("AY47").Value = "'"
# Quote-prefix from the apostrophe trick above adds a style; strip it back
# off so these cells keep the workbook's plain default formatting.
$ws.Range("I47").Style = "Normal"
$ws.Range("N47").Style = "Normal"
$ws.Range("Y47").Style = "Normal"
$ws.Range("AA47").Style = "Normal"
$ws.Range("AT47").Style = "Normal"
$ws.Range("AY47").Style = "Normal"

# Row 48 - Talltita (Poecile montanus)
$ws.Range("A48").Value = 131082526
$ws.Range("B48").Value = 58043
$ws.Range("D48").Value = "NT"
$ws.Range("E48").Value = 103021
$ws.Range("F48").Value = "Talltita"
$ws.Range("G48").Value = "Poecile montanus"
$ws.Range("H48").Value = "(Conrad von Baldenstein, 1827)"
$ws.Range("I48").Value = "'2"
$ws.Range("K48").Value = "adult"
$ws.Range("L48").Value = "'"
$ws.Range("M48").Value = "permanent revir"
$ws.Range("N48").Value = "'"
$ws.Range("P48").Value = "Glupen- nedklassad naturvårdsskog, Ög"
$ws.Range("Q48").Value = 571389
$ws.Range("R48").Value = 6465724
$ws.Range("S48").Value = 10
$ws.Range("T48").Value = "Östergötland"
$ws.Range("U48").Value = "Söderköping"
$ws.Range("V48").Value = "Östergötland"
$ws.Range("W48").Value = "Östra Ryd"
$ws.Range("Y48").Value = "'2026-02-07"
$ws.Range("AA48").Value = "'2026-02-07"
$ws.Range("AC48").Value = "Revirparet"
$ws.Range("AD48").Value = $false
$ws.Range("AE48").Value = $false
$ws.Range("AG48").Value = $false
$ws.Range("AT48").Value = "'"
$ws.Range("AW48").Value = "Steve Daurer"
$ws.Range("AX48").Value = "Steve Daurer"
$ws.Range("AY48").Value = "'"
$ws.Range("I48").Style = "Normal"
$ws.Range("L48").Style = "Normal"
$ws.Range("N48").Style = "Normal"
$ws.Range("Y48").Style = "Normal"
$ws.Range("AA48").Style = "Normal"
$ws.Range("AT48").Style = "Normal"
$ws.Range("AY48").Style = "Normal"
